$wb = $excel.ActiveWorkbook

# --- Service sheet: Brij Mohan Sharma's effort re-categorized from
#     "Operation Management" to "Project Management" ---
$wsService = $wb.Worksheets.Item("Service")
$wsService.Range("C5").Value = "Project Management"

# --- Tourism sheet: new entry for updating the business model ---
$wsTourism = $wb.Worksheets.Item("Tourism")
$wsTourism.Range("A2").Value = "23.10.2016"
$wsTourism.Range("B2").Value = "Arpan Kar"
$wsTourism.Range("C2").Value = "Operation Management"
$wsTourism.Range("D2").Value = 3
$wsTourism.Range("J2").Value = "Updating business model to target April 2016 launch"

# --- Update selections/active cells on a few sheets, and switch the
#     active tab from Service to Cover, matching the user's final
#     navigation state ---
$wsRetail = $wb.Worksheets.Item("Retail")
$wsRetail.Range("A2").Select()

$wsTourism.Range("A3").Select()

$wsTeam = $wb.Worksheets.Item("Team Member")
$wsTeam.Range("J5").Select()

$wsService.Range("C4").Select()

$wsCover = $wb.Worksheets.Item("Cover")
$wsCover.Select()
